$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values for rows 2-6
$ws.Range("D2").Value = 236
$ws.Range("E2").Value = -27
$ws.Range("F2").Value = -27
$ws.Range("G2").Value = -34
$ws.Range("H2").Value = -31
$ws.Range("I2").Value = -31
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 705
$ws.Range("L2").Value = 218
$ws.Range("M2").Value = 487
$ws.Range("N2").Value = 486
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 214
$ws.Range("Q2").Value = -167
$ws.Range("R2").Value = -5
$ws.Range("S2").Value = -25
$ws.Range("T2").Value = 11
$ws.Range("U2").Value = -178
$ws.Range("V2").Value = 95
$ws.Range("W2").Value = -11.44
$ws.Range("X2").Value = -13.2
$ws.Range("Y2").Value = -6.1
$ws.Range("Z2").Value = -3.83
$ws.Range("AA2").Value = 44.85
$ws.Range("AB2").Value = 134.32
$ws.Range("AC2").Value = -71
$ws.Range("AD2").Value = -30.89
$ws.Range("AE2").Value = 1249
$ws.Range("AF2").Value = 1.77
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 42851600
$ws.Range("D3").Value = 239
$ws.Range("E3").Value = -19
$ws.Range("F3").Value = -19
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1154
$ws.Range("L3").Value = 638
$ws.Range("M3").Value = 517
$ws.Range("N3").Value = 509
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 214
$ws.Range("Q3").Value = -553
$ws.Range("R3").Value = 111
$ws.Range("S3").Value = 433
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = -555
$ws.Range("V3").Value = 523
$ws.Range("W3").Value = -7.78
$ws.Range("X3").Value = 5.08
$ws.Range("Y3").Value = 2.36
$ws.Range("Z3").Value = 1.31
$ws.Range("AA3").Value = 123.44
$ws.Range("AB3").Value = 145.14
$ws.Range("AC3").Value = 27
$ws.Range("AD3").Value = 53.39
$ws.Range("AE3").Value = 1308
$ws.Range("AF3").Value = 1.12
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 42851600
$ws.Range("D4").Value = 137
$ws.Range("E4").Value = -44
$ws.Range("F4").Value = -44
$ws.Range("G4").Value = -98
$ws.Range("H4").Value = -95
$ws.Range("I4").Value = -81
$ws.Range("J4").Value = -14
$ws.Range("K4").Value = 1173
$ws.Range("L4").Value = 751
$ws.Range("M4").Value = 422
$ws.Range("N4").Value = 429
$ws.Range("O4").Value = -7
$ws.Range("P4").Value = 214
$ws.Range("Q4").Value = -11
$ws.Range("R4").Value = -50
$ws.Range("S4").Value = 56
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = -11
$ws.Range("V4").Value = 578
$ws.Range("W4").Value = -32.41
$ws.Range("X4").Value = -69.56999999999999
$ws.Range("Y4").Value = -17.31
$ws.Range("Z4").Value = -8.199999999999999
$ws.Range("AA4").Value = 177.72
$ws.Range("AB4").Value = 108.09
$ws.Range("AC4").Value = -190
$ws.Range("AD4").Value = -9.26
$ws.Range("AE4").Value = 1104
$ws.Range("AF4").Value = 1.59
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 42851600
$ws.Range("D5").Value = 539
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = -24
$ws.Range("H5").Value = -30
$ws.Range("I5").Value = -19
$ws.Range("J5").Value = -11
$ws.Range("K5").Value = 881
$ws.Range("L5").Value = 488
$ws.Range("M5").Value = 393
$ws.Range("N5").Value = 414
$ws.Range("O5").Value = -20
$ws.Range("P5").Value = 214
$ws.Range("Q5").Value = 242
$ws.Range("R5").Value = 20
$ws.Range("S5").Value = -226
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 237
$ws.Range("V5").Value = 354
$ws.Range("W5").Value = 0.05
$ws.Range("X5").Value = -5.54
$ws.Range("Y5").Value = -4.42
$ws.Range("Z5").Value = -2.9
$ws.Range("AA5").Value = 124.25
$ws.Range("AB5").Value = 100.67
$ws.Range("AC5").Value = -43
$ws.Range("AD5").Value = -18.95
$ws.Range("AE5").Value = 1063
$ws.Range("AF5").Value = 0.78
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 42851600
$ws.Range("D6").Value = 266
$ws.Range("E6").Value = -27
$ws.Range("F6").Value = -27
$ws.Range("G6").Value = -40
$ws.Range("H6").Value = -52
$ws.Range("I6").Value = -42
$ws.Range("K6").Value = 803
$ws.Range("L6").Value = 459
$ws.Range("M6").Value = 344
$ws.Range("N6").Value = 369
$ws.Range("P6").Value = 214
$ws.Range("Q6").Value = -28
$ws.Range("R6").Value = -6
$ws.Range("S6").Value = -13
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = -29
$ws.Range("V6").Value = 333
$ws.Range("W6").Value = -10.09
$ws.Range("X6").Value = -19.66
$ws.Range("Y6").Value = -10.63
$ws.Range("Z6").Value = -6.2
$ws.Range("AA6").Value = 133.53
$ws.Range("AB6").Value = 80.05
$ws.Range("AC6").Value = -97
$ws.Range("AD6").Value = -6.44
$ws.Range("AE6").Value = 950
$ws.Range("AF6").Value = 0.66
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 42851600

# Clear cells removed entirely in rows 4-6 (AG/AH)
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Clear all data cells D:AI for rows 7-9 (future estimate rows removed)
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
